$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "Test1@34"

# Update row 3
$ws.Range("A3").Value = "tester2"
$ws.Range("B3").Value = "something"

# Remove rows 4-8 (previously had data, now should be empty so dimension shrinks to A1:B3)
$ws.Range("A4:B8").ClearContents()
